$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.1450045
$ws.Range("H2").Value = 14.290009
$ws.Range("I2").Value = 0.8119737125238713
$ws.Range("J2").Value = 0.7990590344890214
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.0003185
$ws.Range("N2").Value = 0.000637
$ws.Range("O2").Value = 0.001593934541086978
$ws.Range("P2").Value = 0.001183922444716212
$ws.Range("Q2").Value = 0.00227568393325
$ws.Range("R2").Value = 0.009102735733
$ws.Range("S2").Value = 0.001294232946846427
$ws.Range("T2").Value = 0.0009460239255848181

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.1450045
$ws.Range("H3").Value = 14.290009
$ws.Range("I3").Value = 0.8119737125238713
$ws.Range("J3").Value = 0.7990590344890214
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.138402
$ws.Range("N3").Value = 0.415206
$ws.Range("O3").Value = 0.6926333700330297
$ws.Range("P3").Value = 0.7716981202210981
$ws.Range("Q3").Value = 0.988882912809
$ws.Range("R3").Value = 5.933297476853999
$ws.Range("S3").Value = 0.5624000888836395
$ws.Range("T3").Value = 0.6166323548608633

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.1450045
$ws.Range("H4").Value = 14.290009
$ws.Range("I4").Value = 0.8119737125238713
$ws.Range("J4").Value = 0.7990590344890214
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.0610995
$ws.Range("N4").Value = 0.122199
$ws.Range("O4").Value = 0.3057726954258833
$ws.Range("P4").Value = 0.2271179573341859
$ws.Range("Q4").Value = 0.43655620244775
$ws.Range("R4").Value = 1.746224809791
$ws.Range("S4").Value = 0.2482793906933854
$ws.Range("T4").Value = 0.1814806557025733

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.284443
$ws.Range("H5").Value = 0.853329
$ws.Range("I5").Value = 0.0323247156403369
$ws.Range("J5").Value = 0.04771587245616726
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.0003185
$ws.Range("N5").Value = 0.000637
$ws.Range("O5").Value = 0.001593934541086978
$ws.Range("P5").Value = 0.001183922444716212
$ws.Range("Q5").Value = 0.0000905950955
$ws.Range("R5").Value = 0.0005435705729999999
$ws.Range("S5").Value = 0.00005152348078994746
$ws.Range("T5").Value = 0.0000564918923700725

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.284443
$ws.Range("H6").Value = 0.853329
$ws.Range("I6").Value = 0.0323247156403369
$ws.Range("J6").Value = 0.04771587245616726
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.138402
$ws.Range("N6").Value = 0.415206
$ws.Range("O6").Value = 0.6926333700330297
$ws.Range("P6").Value = 0.7716981202210981
$ws.Range("Q6").Value = 0.039367480086
$ws.Range("R6").Value = 0.354307320774
$ws.Range("S6").Value = 0.02238917672932593
$ws.Range("T6").Value = 0.03682224907913394

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.284443
$ws.Range("H7").Value = 0.853329
$ws.Range("I7").Value = 0.0323247156403369
$ws.Range("J7").Value = 0.04771587245616726
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.0610995
$ws.Range("N7").Value = 0.122199
$ws.Range("O7").Value = 0.3057726954258833
$ws.Range("P7").Value = 0.2271179573341859
$ws.Range("Q7").Value = 0.0173793250785
$ws.Range("R7").Value = 0.104275950471
$ws.Range("S7").Value = 0.009884015430221021
$ws.Range("T7").Value = 0.01083713148466325

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ntrk1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.370104
$ws.Range("H8").Value = 2.740208
$ws.Range("I8").Value = 0.1557015718357919
$ws.Range("J8").Value = 0.1532250930548114
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.0003185
$ws.Range("N8").Value = 0.000637
$ws.Range("O8").Value = 0.001593934541086978
$ws.Range("P8").Value = 0.001183922444716212
$ws.Range("Q8").Value = 0.000436378124
$ws.Range("R8").Value = 0.001745512496
$ws.Range("S8").Value = 0.000248178113450604
$ws.Range("T8").Value = 0.0001814066267613214

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ntrk1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.370104
$ws.Range("H9").Value = 2.740208
$ws.Range("I9").Value = 0.1557015718357919
$ws.Range("J9").Value = 0.1532250930548114
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.138402
$ws.Range("N9").Value = 0.415206
$ws.Range("O9").Value = 0.6926333700330297
$ws.Range("P9").Value = 0.7716981202210981
$ws.Range("Q9").Value = 0.189625133808
$ws.Range("R9").Value = 1.137750802848
$ws.Range("S9").Value = 0.1078441044200644
$ws.Range("T9").Value = 0.1182435162811008

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntf3"
$ws.Range("C10").Value = "Ntrk1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.370104
$ws.Range("H10").Value = 2.740208
$ws.Range("I10").Value = 0.1557015718357919
$ws.Range("J10").Value = 0.1532250930548114
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.0610995
$ws.Range("N10").Value = 0.122199
$ws.Range("O10").Value = 0.3057726954258833
$ws.Range("P10").Value = 0.2271179573341859
$ws.Range("Q10").Value = 0.08371266934800001
$ws.Range("R10").Value = 0.334850677392
$ws.Range("S10").Value = 0.04760928930227688
$ws.Range("T10").Value = 0.03480017014694931
